# Sync attendance_reports: normalize the "Recorded By" (column G) audit list
# ordering. The recorder-tracking system had been writing this list with a
# stale ordering where the *last* two contributors were transposed; swap
# them back so the most-recent editor consistently appears last, matching
# the canonical ordering used across modules_schedules/assets.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $raw = $cell.Value2

    if ($raw -eq $null) { continue }
    if (-not ($raw -is [string])) { continue }

    # Leave the known-good "backup@backdoor.com, System" pairing untouched.
    # NOTE: use .Equals() (ordinal, case-sensitive) for comparisons here —
    # the `-eq`/`-ne` operators in this host fold case, which would wrongly
    # treat "...System, system" / "...system, System" as identical.
    if ($raw.Equals("backup@backdoor.com, System")) { continue }

    $parts = $raw -split ", "
    if ($parts.Count -lt 2) { continue }

    $lastIdx = $parts.Count - 1
    $tmp = $parts[$lastIdx]
    $parts[$lastIdx] = $parts[$lastIdx - 1]
    $parts[$lastIdx - 1] = $tmp

    $newVal = [string]::Join(", ", $parts)
    if (-not $newVal.Equals($raw)) {
        $cell.Value = $newVal
    }
}
